$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44194
$ws.Range("I2").Value = 'Extra'
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 3500
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = 3500
$ws.Range("N2").Value = '$/unidad'
$ws.Range("O2").Value = 'Región de O''Higgins'
$ws.Range("P2").Value = 3500

# Row 3
$ws.Range("D3").Value = 44194
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 3000
$ws.Range("P3").Value = 3000

# Row 4
$ws.Range("D4").Value = 44167
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = 5000
$ws.Range("P4").Value = 5000

# Row 5
$ws.Range("D5").Value = 44167
$ws.Range("I5").Value = 'Segunda'
$ws.Range("J5").Value = 560
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 3000
$ws.Range("N5").Value = '$/unidad'
$ws.Range("O5").Value = 'Región de O''Higgins'
$ws.Range("P5").Value = 3000

# Row 6
$ws.Range("D6").Value = 44167
$ws.Range("I6").Value = 'Tercera'
$ws.Range("J6").Value = 450
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 2000
$ws.Range("P6").Value = 2000

# Row 7
$ws.Range("D7").Value = 44495
$ws.Range("K7").Value = 800
$ws.Range("L7").Value = 800
$ws.Range("M7").Value = 800
$ws.Range("N7").Value = '$/kilo (volumen en unidades)'
$ws.Range("O7").Value = 'Perú'
$ws.Range("P7").Value = 800

# Row 8
$ws.Range("D8").Value = 44510
$ws.Range("J8").Value = 250

# Row 9
$ws.Range("D9").Value = 44477
$ws.Range("J9").Value = 80

# Row 10
$ws.Range("D10").Value = 44497

# Row 11
$ws.Range("D11").Value = 44488
$ws.Range("J11").Value = 150

# Row 12
$ws.Range("D12").Value = 44217
$ws.Range("I12").Value = 'Extra'
$ws.Range("K12").Value = 2500
$ws.Range("L12").Value = 2500
$ws.Range("M12").Value = 2500
$ws.Range("P12").Value = 2500

# Row 13
$ws.Range("D13").Value = 44217
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 280
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = 2000
$ws.Range("P13").Value = 2000

# Row 14
$ws.Range("D14").Value = 44491
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 150
$ws.Range("K14").Value = 800
$ws.Range("L14").Value = 800
$ws.Range("M14").Value = 800
$ws.Range("N14").Value = '$/kilo (volumen en unidades)'
$ws.Range("O14").Value = 'Perú'
$ws.Range("P14").Value = 800

# Row 15
$ws.Range("D15").Value = 44483
$ws.Range("J15").Value = 120

# Row 16
$ws.Range("D16").Value = 44312
$ws.Range("J16").Value = 180
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = 2500
$ws.Range("N16").Value = '$/unidad'
$ws.Range("P16").Value = 2500

# Row 18
$ws.Range("D18").Value = 44504
$ws.Range("H18").Value = 'Sin especificar'
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 800
$ws.Range("L18").Value = 800
$ws.Range("M18").Value = 800
$ws.Range("N18").Value = '$/kilo (volumen en unidades)'
$ws.Range("O18").Value = 'Perú'
$ws.Range("P18").Value = 800

# Row 19
$ws.Range("I19").Value = 'Extra'
$ws.Range("J19").Value = 340
$ws.Range("K19").Value = 2500
$ws.Range("L19").Value = 2500
$ws.Range("M19").Value = 2500
$ws.Range("P19").Value = 2500

# Row 20
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = 2000
$ws.Range("P20").Value = 2000

# Row 21
$ws.Range("I21").Value = 'Segunda'
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 1500
$ws.Range("M21").Value = 1500
$ws.Range("P21").Value = 1500

# Row 22
$ws.Range("D22").Value = 44223
$ws.Range("H22").Value = 'Americana O Klondike'
$ws.Range("I22").Value = 'Tercera'
$ws.Range("J22").Value = 160
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 1000
$ws.Range("O22").Value = 'Región de O''Higgins'
$ws.Range("P22").Value = 1000
